$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper re-ran and the two matches that shared identical kickoff
# timestamps came back in the opposite order; for each such pair the match
# detail columns (F:V) swap between the two rows while the index/meta
# columns (A:E) stay put.
$swapPairs = @(
    @(32, 33),
    @(35, 36),
    @(37, 38),
    @(60, 61),
    @(64, 65),
    @(66, 67),
    @(80, 81),
    @(84, 85),
    @(96, 97),
    @(102, 103),
    @(110, 111),
    @(125, 126)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("F$r1" + ":V$r1")
    $range2 = $ws.Range("F$r2" + ":V$r2")
    $vals1 = $range1.Value()
    $vals2 = $range2.Value()
    $range1.Value = $vals2
    $range2.Value = $vals1
}

# Two new matches were appended at the bottom of the sheet. Clone the
# formatting of the last existing data row (169) onto the two new rows,
# then overwrite with the new match data.
$ws.Range("A169:V169").Copy()
$ws.Range("A170:V170").PasteSpecial()
$ws.Range("A169:V169").Copy()
$ws.Range("A171:V171").PasteSpecial()

$ws.Range("A170").Value = 169
$ws.Range("F170").Value = "Central Cordoba"
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = "Estudiantes L.P."
$ws.Range("I170").Value = 1
$ws.Range("J170").Value = 3.15
$ws.Range("K170").Value = "07/11/2023 06:42"
$ws.Range("L170").Value = 3.36
$ws.Range("M170").Value = "10/11/2023 22:49"
$ws.Range("N170").Value = 3
$ws.Range("O170").Value = "07/11/2023 06:42"
$ws.Range("P170").Value = 3.06
$ws.Range("Q170").Value = "10/11/2023 22:54"
$ws.Range("R170").Value = 2.55
$ws.Range("S170").Value = "07/11/2023 06:42"
$ws.Range("T170").Value = 2.44
$ws.Range("U170").Value = "10/11/2023 22:49"
$ws.Range("V170").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/central-cordoba-santiago-del-estero-estudiantes-l-p/Q3k9WAuP/"
$ws.Range("E170").Value = 45240.95833333334

$ws.Range("A171").Value = 170
$ws.Range("F171").Value = "Gimnasia L.P."
$ws.Range("G171").Value = 1
$ws.Range("H171").Value = "Atl. Tucuman"
$ws.Range("I171").Value = 2
$ws.Range("J171").Value = 2.21
$ws.Range("K171").Value = "07/11/2023 06:42"
$ws.Range("L171").Value = 2.31
$ws.Range("M171").Value = "11/11/2023 00:58"
$ws.Range("N171").Value = 3.13
$ws.Range("O171").Value = "07/11/2023 06:42"
$ws.Range("P171").Value = 3.06
$ws.Range("Q171").Value = "11/11/2023 00:58"
$ws.Range("R171").Value = 3.68
$ws.Range("S171").Value = "07/11/2023 06:42"
$ws.Range("T171").Value = 3.65
$ws.Range("U171").Value = "11/11/2023 00:58"
$ws.Range("V171").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/gimnasia-l-p-atl-tucuman/CKTmweds/"
$ws.Range("E171").Value = 45241.04166666666

# Re-apply the bold / bordered / centred look that column A uses for the
# index, since writing a new value resets a cell's font/border/alignment.
foreach ($r in @(170, 171)) {
    $cell = $ws.Range("A$r")
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}
